$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = -7.328999999999999
$ws.Range("D27").Value = -8.794
$ws.Range("D32").Value = -8.010000000000002
$ws.Range("D36").Value = -8.051
$ws.Range("D38").Value = -7.662000000000001
$ws.Range("D46").Value = -8.122
$ws.Range("D54").Value = -8.558
$ws.Range("D55").Value = -8.129000000000001
$ws.Range("D56").Value = -8.403
$ws.Range("D67").Value = -7.281000000000001
$ws.Range("D69").Value = -7.321000000000001
$ws.Range("D72").Value = -7.434
$ws.Range("D83").Value = -8.047000000000001
$ws.Range("D86").Value = -8.118999999999998
$ws.Range("D91").Value = -6.970000000000001
$ws.Range("D93").Value = -7.568
$ws.Range("D99").Value = -8.019
